$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 36-43 (table shrank from 41 tickers to 33 tickers in the
# "support Zone" column; dimension goes from A1:F43 to A1:F35).
# Delete bottom-up so row numbers of not-yet-deleted rows stay stable.
$ws.Range("A43:F43").EntireRow.Delete()
$ws.Range("A42:F42").EntireRow.Delete()
$ws.Range("A41:F41").EntireRow.Delete()
$ws.Range("A40:F40").EntireRow.Delete()
$ws.Range("A39:F39").EntireRow.Delete()
$ws.Range("A38:F38").EntireRow.Delete()
$ws.Range("A37:F37").EntireRow.Delete()
$ws.Range("A36:F36").EntireRow.Delete()

# Update the watch-list ticker values that changed across columns B, C, D, E, F
$ws.Range("B2").Value = "NSE:DEVYANI"
$ws.Range("C2").Value = "NSE:ADSL"
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = "NSE:BAJAJ-AUTO"
$ws.Range("F2").Value = $null
$ws.Range("B3").Value = "NSE:EMMBI"
$ws.Range("C3").Value = "NSE:AKSHARCHEM"
$ws.Range("E3").Value = "NSE:BIOCON"
$ws.Range("B4").Value = "NSE:NH"
$ws.Range("C4").Value = "NSE:APCL"
$ws.Range("E4").Value = "NSE:BOSCHLTD"
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = "NSE:ASIANPAINT"
$ws.Range("E5").Value = "NSE:CDSL"
$ws.Range("C6").Value = "NSE:ASMS"
$ws.Range("E6").Value = "NSE:CGPOWER"
$ws.Range("C7").Value = "NSE:AURIONPRO"
$ws.Range("E7").Value = "NSE:COLPAL"
$ws.Range("C8").Value = "NSE:BANSWRAS"
$ws.Range("E8").Value = "NSE:CYIENT"
$ws.Range("C9").Value = "NSE:BATAINDIA"
$ws.Range("E9").Value = "NSE:DALBHARAT"
$ws.Range("C10").Value = "NSE:BEDMUTHA"
$ws.Range("E10").Value = "NSE:GNFC"
$ws.Range("C11").Value = "NSE:BHANDARI"
$ws.Range("E11").Value = "NSE:HFCL"
$ws.Range("C12").Value = "NSE:BPL"
$ws.Range("E12").Value = "NSE:ICICIPRULI"
$ws.Range("C13").Value = "NSE:CENTENKA"
$ws.Range("E13").Value = "NSE:IEX"
$ws.Range("C14").Value = "NSE:CENTRUM"
$ws.Range("E14").Value = "NSE:JSWENERGY"
$ws.Range("C15").Value = "NSE:ENIL"
$ws.Range("E15").Value = "NSE:LALPATHLAB"
$ws.Range("C16").Value = "NSE:EXPLEOSOL"
$ws.Range("E16").Value = "NSE:LTIM"
$ws.Range("C17").Value = "NSE:GENUSPAPER"
$ws.Range("E17").Value = "NSE:LTTS"
$ws.Range("C18").Value = "NSE:GOKEX"
$ws.Range("E18").Value = "NSE:MARUTI"
$ws.Range("C19").Value = "NSE:GULFPETRO"
$ws.Range("E19").Value = "NSE:MCX"
$ws.Range("C20").Value = "NSE:HLVLTD"
$ws.Range("E20").Value = "NSE:NYKAA"
$ws.Range("C21").Value = "NSE:HMAAGRO"
$ws.Range("E21").Value = "NSE:PIIND"
$ws.Range("C22").Value = "NSE:INDORAMA"
$ws.Range("C23").Value = "NSE:INDOWIND"
$ws.Range("C24").Value = "NSE:ISGEC"
$ws.Range("C25").Value = "NSE:JKPAPER"
$ws.Range("C26").Value = "NSE:MOL"
$ws.Range("C27").Value = "NSE:NDLVENTURE"
$ws.Range("C28").Value = "NSE:NUCLEUS"
$ws.Range("C29").Value = "NSE:OIL"
$ws.Range("C30").Value = "NSE:PLAZACABLE"
$ws.Range("C31").Value = "NSE:PNCINFRA"
$ws.Range("C32").Value = "NSE:QUICKHEAL"
$ws.Range("C33").Value = "NSE:RBL"
$ws.Range("C34").Value = "NSE:RCF"
$ws.Range("C35").Value = "NSE:RUBYMILLS"